$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($rowIndex, $newText) {
    $cell = $tbl.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

# Simple single-value replacements
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "24"
Set-CellText 8 "0.00001"
Set-CellText 12 "0.00103"

# Rows that collapse a multi-run / tab-separated breakdown into a single value
Set-CellText 44 "100"
Set-CellText 45 "0"
Set-CellText 46 "60"
